# Apply the "add 2022-Q4 data" edit:
#  1. Insert a new summary row for 2022-Q4 at the top of the "总计" sheet's
#     data table (pushing the existing quarters down by one row).
#  2. Insert a brand-new worksheet named "2022-Q4" (positioned right after
#     "总计" and before "2022-Q3") holding the per-fund breakdown for the
#     new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert the 2022-Q4 row at the top of the data
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
$summary.Rows.Item(2).ClearFormats()

# Match the look of the other index-column cells (bold/centered/bordered).
$summary.Range("A4").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 9
$summary.Range("D2").Value = 1.28

# Re-sequence the 0-based index column now that there are 8 data rows.
for ($i = 0; $i -le 7; $i++) {
    $summary.Cells.Item($i + 2, 1).Value = $i
}

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet with the per-fund detail table
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Header row, copied (with formatting) from the 2022-Q3 sheet so the new
# sheet matches the look of the rest of the workbook.
$q3.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4104)
$excel.CutCopyMode = $false

$q3.Range("A2").Copy()
$q4.Range("A2:A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Cells B:G hold numeric-looking text (fund codes / percentages) stored as
# plain text in the source data -- force a Text format so values such as
# "009010" or "27.18" are not silently coerced into numbers.
$q4.Range("B2:G10").NumberFormat = "@"

$rows = @(
    @(0, "009010", "华夏兴阳一年持有期混合",     "27.18", "90.65", "3.24", "0.8806", 6),
    @(1, "009223", "宝盈现代服务业混合A",         "3.46",  "88.59", "5.72", "0.1979", 6),
    @(2, "013859", "宝盈品质甄选混合A",           "1.83",  "88.38", "5.57", "0.1019", 6),
    @(3, "008303", "宝盈龙头优选股票A",           "0.65",  "88.21", "5.60", "0.0364", 6),
    @(4, "009224", "宝盈现代服务业混合C",         "0.41",  "88.59", "5.72", "0.0235", 6),
    @(5, "006675", "宝盈品牌消费股票A",           "0.18",  "87.69", "5.71", "0.0103", 5),
    @(6, "008304", "宝盈龙头优选股票C",           "0.18",  "88.21", "5.60", "0.0101", 6),
    @(7, "013860", "宝盈品质甄选混合C",           "0.18",  "88.38", "5.57", "0.0100", 6),
    @(8, "006676", "宝盈品牌消费股票C",           "0.15",  "87.69", "5.71", "0.0086", 5)
)

$r = 2
foreach ($row in $rows) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 7).Value = $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
